# Insert a new "Номер справки" (reference number) column at the very
# front of the patient table and fill it in with a reference number for
# each patient row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing A:J columns one slot to the right, opening up a
# fresh column A for the new field.
$ws.Range("A:A").Insert()

# Give the new header cell the same bold/centered look as the rest of
# the header row.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# The old "ИНН"/"Паспорт" header cells (now H1/I1) drop their special
# text-number-format style and pick up the plain bold-centered style
# used by the rest of the header row.
$ws.Range("B1").Copy()
$ws.Range("H1:I1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Header text
$ws.Range("A1").Value = "Номер справки"

# Reference numbers for the 4 patient rows
$ws.Range("A2").Value = 3334
$ws.Range("A3").Value = 3335
$ws.Range("A4").Value = 3336
$ws.Range("A5").Value = 3337

# Column width for the new column
$ws.Range("A:A").ColumnWidth = 18.5

# Re-select cell to match final cursor position recorded in the workbook
$ws.Range("H13").Select()
